$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.424.90'
$ws.Range('E2').Value = '  -3.90%  '
$ws.Range('D3').Value = '1.770.79'
$ws.Range('E3').Value = '  -2.99%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('E5').Value = '  +0.14%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '306.49'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.88%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4309'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.35%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3665'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.54%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07241'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.54%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8515'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.18%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.28'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.54%  '
$ws.Range('D12').Value = '1.790.33'
$ws.Range('E12').Value = '  -1.46%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.444'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.48%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.243'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.52%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06836'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.47%  '
$ws.Range('E16').Value = '  +0.48%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '79.46'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.50%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008690'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.39%  '
$ws.Range('E19').Value = '  +0.17%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.01'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.33%  '
$ws.Range('D21').Value = '26.429.40'
$ws.Range('E21').Value = '  -3.76%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.121'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.27'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.77%  '
$ws.Range('D24').Value = '2.014.13'
$ws.Range('E24').Value = '  -1.13%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.03'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.98%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.850'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -7.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.15'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.91%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.083'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.33%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '114.79'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.58%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.725'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.56%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08939'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.73%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7249'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.22%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.114'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.55%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.325'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.57%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.748'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -7.62%  '
$ws.Range('E36').Value = '  +0.16%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.081'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.40%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05161'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.26%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01895'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.25%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.4931'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.71%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1605'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.00%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.535'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -9.33%  '
$ws.Range('E43').Value = '  -3.26%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.032'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.65%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '104.86'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.32%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.19'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.30%  '
$ws.Range('E47').Value = '  +0.18%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.06200'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.12%  '
$ws.Range('B49').Value = 'Decentraland'
$ws.Range('C49').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4484'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.02%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.585'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.82%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.741'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.21%  '
